# "Device" -> "Apparatus" rename across the IEEE_30Bus workbook.
#
# The only sheet affected is the one literally named "Device": its tab
# name and three of its header/label cells mention "device(s)" and need
# to become "Apparatus"/"apparatuses". All other sheets/cells are
# untouched. Finally, make the renamed sheet the active tab with A2
# selected (matching the author's final interactive state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Rename the worksheet tab.
$ws.Name = "Apparatus"

# Update the textual "device" -> "apparatus" labels on the sheet.
# (Ordered to match the string-table layout produced by the original
# edit: type/parameters labels first, the summary sentence last.)
$ws.Range("B2").Value = "Apparatus type"
$ws.Range("C2").Value = "Apparatus parameters"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Make the renamed sheet the active one, with A2 selected.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
